# Weekly update for "Hortaliza, Comercializadora del Agro de Limarí - Tomate".
#
# A new week of price data (Fecha = 45077, i.e. 2023-05-31) is inserted at the
# top of the data block (rows 858-863), pushing every subsequent row down by
# six rows (the block's weekly group size: Larga vida/Semiduro x
# Primera/Segunda/Tercera). The sheet's used range grows from A1:R958 to
# A1:R964 as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows right before the current row 858, shifting the old
# rows 858-958 down to 864-964 (formatting, incl. the date number format on
# column D, comes along for the ride the way a real Excel row-insert works).
$ws.Rows("858:863").Insert()

# The new week's six records (same market/region/category for this sheet;
# only Variedad/Calidad/Volumen/Precios differ row to row).
$newRows = @(
  @{r=858; H="Larga vida"; I="Primera"; J=1100; K=11000; L=12000; M=11500; P=639},
  @{r=859; H="Larga vida"; I="Segunda"; J=800;  K=9000;  L=10000; M=9500;  P=528},
  @{r=860; H="Larga vida"; I="Tercera"; J=500;  K=7000;  L=8000;  M=7500;  P=417},
  @{r=861; H="Semiduro";   I="Primera"; J=700;  K=8000;  L=9000;  M=8500;  P=472},
  @{r=862; H="Semiduro";   I="Segunda"; J=500;  K=6000;  L=7000;  M=6500;  P=361},
  @{r=863; H="Semiduro";   I="Tercera"; J=400;  K=4000;  L=5000;  M=4500;  P=250}
)

foreach ($row in $newRows) {
    $rn = $row.r
    $ws.Range("A$rn").Value = 2
    $ws.Range("B$rn").Value = "Comercializadora del Agro de Limar" + [char]0x00ED
    $ws.Range("C$rn").Value = "Coquimbo"
    $ws.Range("D$rn").Value = 45077
    $ws.Range("E$rn").Value = 4
    $ws.Range("F$rn").Value = 100112020
    $ws.Range("G$rn").Value = "Tomate"
    $ws.Range("H$rn").Value = $row.H
    $ws.Range("I$rn").Value = $row.I
    $ws.Range("J$rn").Value = $row.J
    $ws.Range("K$rn").Value = $row.K
    $ws.Range("L$rn").Value = $row.L
    $ws.Range("M$rn").Value = $row.M
    $ws.Range("N$rn").Value = "$/bandeja 18 kilos"
    $ws.Range("O$rn").Value = "Provincia de Limar" + [char]0x00ED
    $ws.Range("P$rn").Value = $row.P
    $ws.Range("Q$rn").Value = 18
    $ws.Range("R$rn").Value = "Hortaliza"
}
